$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28").Value = 3979.85302734375
$ws.Range("E28").Value = [double]"5.4060328436023042e-10"
$ws.Range("C29").Value = 5969.779296875
$ws.Range("E29").Value = [double]"7.1604633333777201e-10"
$ws.Range("C30").Value = 11939.8681640625
$ws.Range("E30").Value = [double]"3.6590239727019025e-09"
$ws.Range("C31").Value = 23879.921875
$ws.Range("E31").Value = [double]"7.6143660265870494e-09"
$ws.Range("C32").Value = 45768.3125
$ws.Range("E32").Value = [double]"1.9715924892693693e-08"
$ws.Range("C33").Value = 77610.2578125
$ws.Range("E33").Value = [double]"2.4852688795817812e-08"
$ws.Range("C34").Value = 117786.9765625
$ws.Range("E34").Value = [double]"2.8075906044477961e-08"
$ws.Range("C35").Value = 118327.9453125
$ws.Range("E35").Value = [double]"1.6240404576706169e-08"
$ws.Range("C36").Value = 56474.09765625
$ws.Range("E36").Value = [double]"1.5849668244527493e-08"
$ws.Range("C37").Value = 26512.328125
$ws.Range("E37").Value = [double]"8.1652329342318808e-09"
$ws.Range("C38").Value = 19543.951171875
$ws.Range("E38").Value = [double]"1.0552954066156417e-08"
$ws.Range("C39").Value = 10862.0146484375
$ws.Range("E39").Value = [double]"6.1255596150999736e-09"
$ws.Range("C40").Value = 3729.3203125
$ws.Range("E40").Value = [double]"2.5664343983322624e-09"
$ws.Range("C41").Value = 1345.5067138671875
$ws.Range("E41").Value = [double]"1.7296009069411866e-09"
$ws.Range("C42").Value = 3192.2021484375
$ws.Range("E42").Value = [double]"2.9333384610197299e-09"
$ws.Range("C43").Value = 25.691064834594727
$ws.Range("E43").Value = [double]"3.3238387736611585e-11"
$ws.Range("C44").Value = 254.52424621582031
$ws.Range("E44").Value = [double]"3.1961480684294941e-10"
$ws.Range("C45").Value = 178.62748718261719
$ws.Range("E45").Value = [double]"4.9158421777661943e-10"
$ws.Range("C46").Value = 881.1624755859375
$ws.Range("E46").Value = [double]"2.4023865119460197e-09"
$ws.Range("C76").Value = 29103.9921875
$ws.Range("E76").Value = [double]"4.3031484153743804e-08"
$ws.Range("C77").Value = 50032.640625
$ws.Range("E77").Value = [double]"6.0268301460553175e-09"
$ws.Range("C78").Value = 55395.890625
$ws.Range("E78").Value = [double]"5.5437436863314815e-09"
$ws.Range("C79").Value = 71883.3828125
$ws.Range("E79").Value = [double]"6.3522169746477175e-09"
$ws.Range("C80").Value = 128631.1328125
$ws.Range("E80").Value = [double]"2.9041949289876356e-08"
$ws.Range("C81").Value = 223461.390625
$ws.Range("E81").Value = [double]"5.2494922186951953e-08"
$ws.Range("C82").Value = 140919.203125
$ws.Range("E82").Value = [double]"4.4723563519255549e-08"
$ws.Range("C83").Value = 76175.546875
$ws.Range("E83").Value = [double]"1.7971478527556428e-08"
$ws.Range("C84").Value = 61683.2578125
$ws.Range("E84").Value = [double]"1.0832227559376406e-08"
$ws.Range("C85").Value = 65187.88671875
$ws.Range("E85").Value = [double]"6.5915926050763574e-09"
$ws.Range("C86").Value = 33684.73828125
$ws.Range("E86").Value = [double]"6.9649495060275513e-09"
$ws.Range("C87").Value = 33559.91015625
$ws.Range("E87").Value = [double]"7.6147461669506811e-09"
$ws.Range("C88").Value = 17296.1953125
$ws.Range("E88").Value = [double]"6.8805987574194205e-09"
$ws.Range("C89").Value = 17406.576171875
$ws.Range("E89").Value = [double]"7.2320713861984132e-09"
$ws.Range("C90").Value = 8982.533203125
$ws.Range("E90").Value = [double]"4.554212118534906e-09"
$ws.Range("C91").Value = 1941.474365234375
$ws.Range("E91").Value = [double]"1.8386778766199541e-09"
$ws.Range("C92").Value = 1389.294921875
$ws.Range("E92").Value = [double]"9.4054630750406432e-10"
$ws.Range("C93").Value = 251.79917907714844
$ws.Range("E93").Value = [double]"2.4000823550807127e-10"
$ws.Range("C95").Value = 1.5880948305130005
$ws.Range("E95").Value = [double]"3.2198824769652523e-12"
$ws.Range("C96").Value = 575.76568603515625
$ws.Range("E96").Value = [double]"1.1565025603843537e-09"
$ws.Range("C205").Value = 42157.80029296875
$ws.Range("E205").Value = [double]"4.3215464984314167e-08"
$ws.Range("C206").Value = 72246.0390625
$ws.Range("E206").Value = [double]"6.0336073914868393e-09"
$ws.Range("C207").Value = 79990.453125
$ws.Range("E207").Value = [double]"5.5499778106593567e-09"
$ws.Range("C208").Value = 103798.0234375
$ws.Range("E208").Value = [double]"6.3593592614097361e-09"
$ws.Range("C209").Value = 185740.53125
$ws.Range("E209").Value = [double]"2.9074605834011891e-08"
$ws.Range("C210").Value = 322673.34375
$ws.Range("E210").Value = [double]"5.2553950524725224e-08"
$ws.Range("C211").Value = 203484.234375
$ws.Range("E211").Value = [double]"4.477385218137897e-08"
$ws.Range("C212").Value = 109995.8203125
$ws.Range("E212").Value = [double]"1.7991686362961445e-08"
$ws.Range("C213").Value = 89069.2734375
$ws.Range("E213").Value = [double]"1.0844408038224174e-08"
$ws.Range("C214").Value = 94129.8828125
$ws.Range("E214").Value = [double]"6.5990053421671746e-09"
$ws.Range("C215").Value = 48507.75
$ws.Range("E215").Value = [double]"6.9538192981610791e-09"
$ws.Range("C216").Value = 68175.4453125
$ws.Range("E216").Value = [double]"1.0724822807617329e-08"
$ws.Range("C217").Value = 35136.44140625
$ws.Range("E217").Value = [double]"9.6908285840413555e-09"
$ws.Range("C218").Value = 35360.6796875
$ws.Range("E218").Value = [double]"1.0185853938082801e-08"
$ws.Range("C219").Value = 18247.61328125
$ws.Range("E219").Value = [double]"6.4142811062595229e-09"
$ws.Range("C220").Value = 3944.017822265625
$ws.Range("E220").Value = [double]"2.5896456090634956e-09"
$ws.Range("C221").Value = 2822.2900390625
$ws.Range("E221").Value = [double]"1.3246919117193556e-09"
$ws.Range("C222").Value = 511.51873779296875
$ws.Range("E222").Value = [double]"3.3803435051105168e-10"
$ws.Range("C224").Value = 3.2261433601379395
$ws.Range("E224").Value = [double]"4.5349735251776213e-12"
$ws.Range("C225").Value = 1169.64208984375
$ws.Range("E225").Value = [double]"1.6288508319917128e-09"
